$wb = $excel.ActiveWorkbook

# --- Update data values on the "SoCDTtiNTY-psgr" sheet ---
$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# Row 2 (LDVs): B2 was a formula (=0.076+(0.076-0.0725)); replace with plain value 0.076
$wsPsgr.Range("B2").Value = 0.076
# D2 was 0.0735; update to 0.076
$wsPsgr.Range("D2").Value = 0.076

# Row 5 (aircraft): B5 and E5 were 0.01; update to 0.029
$wsPsgr.Range("B5").Value = 0.029
$wsPsgr.Range("E5").Value = 0.029

# --- Switch the active/selected sheet from "SoCDTtiNTY-psgr" to "About" ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()

Write-Host "done"
